$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8322162453630995
$ws.Range("C2").Value = 0.283447714744625
$ws.Range("D2").Value = 0.07891491415117002
$ws.Range("E2").Value = 0.09666655582798711
$ws.Range("G2").Value = 0.9509109816039825
$ws.Range("H2").Value = 0.8875359085953392
$ws.Range("I2").Value = 0.6726832726202687
$ws.Range("M2").Value = 0.3740316256657437

# Row 3
$ws.Range("B3").Value = 0.7367777205278117
$ws.Range("C3").Value = 0.2473495584821421
$ws.Range("D3").Value = 0.07150616303459856
$ws.Range("E3").Value = 0.09140658319989825
$ws.Range("G3").Value = 0.9125796724722477
$ws.Range("H3").Value = 0.8759039180849584
$ws.Range("I3").Value = 0.664337401815402
$ws.Range("M3").Value = 0.3357718321900478

# Row 4
$ws.Range("B4").Value = 0.6783278838970546
$ws.Range("C4").Value = 0.2251976797865041
$ws.Range("D4").Value = 0.06699813706057967
$ws.Range("E4").Value = 0.0882547871705377
$ws.Range("G4").Value = 0.8898506950751823
$ws.Range("H4").Value = 0.869375389207562
$ws.Range("I4").Value = 0.6597338485623325
$ws.Range("M4").Value = 0.3124341468380578

# Row 5
$ws.Range("B5").Value = 0.6545465186071908
$ws.Range("C5").Value = 0.2161732952422994
$ws.Range("D5").Value = 0.06517122820055476
$ws.Range("E5").Value = 0.08698967664106405
$ws.Range("G5").Value = 0.8807889754599501
$ws.Range("H5").Value = 0.8668683050122468
$ws.Range("I5").Value = 0.6579878026193384
$ws.Range("M5").Value = 0.3029619093792988

# Row 6
$ws.Range("B6").Value = 0.6505998943505062
$ws.Range("C6").Value = 0.2146749519035041
$ws.Range("D6").Value = 0.06486847979728338
$ws.Range("E6").Value = 0.08678076114250288
$ws.Range("G6").Value = 0.8792963277452657
$ws.Range("H6").Value = 0.8664612398067106
$ws.Range("I6").Value = 0.6577056909213752
$ws.Range("M6").Value = 0.3013913311597989

# Row 7
$ws.Range("B7").Value = 0.6780070088923367
$ws.Range("C7").Value = 0.2250759635089139
$ws.Range("D7").Value = 0.06697345784145625
$ws.Range("E7").Value = 0.08823764781488919
$ws.Range("G7").Value = 0.8897276766913365
$ws.Range("H7").Value = 0.86934095816423
$ws.Range("I7").Value = 0.6597097760293096
$ws.Range("M7").Value = 0.3123062476934706

# Row 8
$ws.Range("B8").Value = 0.7992777598744283
$ws.Range("C8").Value = 0.2709981745272501
$ws.Range("D8").Value = 0.0763517754931371
$ws.Range("E8").Value = 0.09483657477851892
$ws.Range("G8").Value = 0.937525534516169
$ws.Range("H8").Value = 0.8833972807814234
$ws.Range("I8").Value = 0.6696968894991286
$ws.Range("M8").Value = 0.3608072765829817

# Row 9
$ws.Range("B9").Value = 1.038303153122229
$ws.Range("C9").Value = 0.3611767906941168
$ws.Range("D9").Value = 0.09507579579374692
$ws.Range("E9").Value = 0.1084092019609599
$ws.Range("G9").Value = 1.03776656085202
$ws.Range("H9").Value = 0.9158750345332578
$ws.Range("I9").Value = 0.693462552073143
$ws.Range("M9").Value = 0.4571746224770976

# Row 10
$ws.Range("B10").Value = 1.214709963447433
$ws.Range("C10").Value = 0.4275501252002982
$ws.Range("D10").Value = 0.1090480186368978
$ws.Range("E10").Value = 0.1187875089191266
$ws.Range("G10").Value = 1.115541907408584
$ws.Range("H10").Value = 0.9427978532500845
$ws.Range("I10").Value = 0.7135414319435966
$ws.Range("M10").Value = 0.5287984694164578

# Row 11
$ws.Range("B11").Value = 1.295148138574916
$ws.Range("C11").Value = 0.4577807403847487
$ws.Range("D11").Value = 0.1154540042637109
$ws.Range("E11").Value = 0.1236017976820136
$ws.Range("G11").Value = 1.151855455529073
$ws.Range("H11").Value = 0.955725084269659
$ws.Range("I11").Value = 0.7232595449898866
$ws.Range("M11").Value = 0.5615735099464985

# Row 12
$ws.Range("B12").Value = 1.325636027449832
$ws.Range("C12").Value = 0.4692342873635766
$ws.Range("D12").Value = 0.1178871689917003
$ws.Range("E12").Value = 0.1254385993121758
$ws.Range("G12").Value = 1.165743339298928
$ws.Range("H12").Value = 0.9607191846277487
$ws.Range("I12").Value = 0.7270247294713386
$ws.Range("M12").Value = 0.5740132236234814

# Row 13
$ws.Range("B13").Value = 1.3190686764612
$ws.Range("C13").Value = 0.466767291297856
$ws.Range("D13").Value = 0.1173628134311997
$ws.Range("E13").Value = 0.1250423948814401
$ws.Range("G13").Value = 1.162746213150228
$ws.Range("H13").Value = 0.9596392027553975
$ws.Range("I13").Value = 0.7262100243572718
$ws.Range("M13").Value = 0.5713328280449161

# Row 14
$ws.Range("B14").Value = 1.29765583801742
$ws.Range("C14").Value = 0.4587229098121384
$ws.Range("D14").Value = 0.1156540337389629
$ws.Range("E14").Value = 0.1237526352219689
$ws.Range("G14").Value = 1.152995265140675
$ws.Range("H14").Value = 0.9561339648821274
$ws.Range("I14").Value = 0.7235675959237113
$ws.Range("M14").Value = 0.5625963567103724

# Row 15
$ws.Range("B15").Value = 1.284543475740804
$ws.Range("C15").Value = 0.4537962745063169
$ws.Range("D15").Value = 0.1146083201458623
$ws.Range("E15").Value = 0.122964419070378
$ws.Range("G15").Value = 1.147040410237736
$ws.Range("H15").Value = 0.9539998090625375
$ws.Range("I15").Value = 0.7219601549549139
$ws.Range("M15").Value = 0.5572487525179213

# Row 16
$ws.Range("B16").Value = 1.209456992806963
$ws.Range("C16").Value = 0.425575274916639
$ws.Range("D16").Value = 0.1086303923825369
$ws.Range("E16").Value = 0.1184747871859102
$ws.Range("G16").Value = 1.113187726165705
$ws.Range("H16").Value = 0.9419668008850692
$ws.Range("I16").Value = 0.71291818261858
$ws.Range("M16").Value = 0.5266604919701052

# Row 17
$ws.Range("B17").Value = 1.163442761324688
$ws.Range("C17").Value = 0.4082724762393468
$ws.Range("D17").Value = 0.1049760441353413
$ws.Range("E17").Value = 0.1157446413582832
$ws.Range("G17").Value = 1.092661008481798
$ws.Range("H17").Value = 0.9347598376499775
$ws.Range("I17").Value = 0.7075216501488626
$ws.Range("M17").Value = 0.5079455541640812

# Row 18
$ws.Range("B18").Value = 1.136994434124688
$ws.Range("C18").Value = 0.3983238184047764
$ws.Range("D18").Value = 0.102878850779291
$ws.Range("E18").Value = 0.1141830803094592
$ws.Range("G18").Value = 1.080942262848509
$ws.Range("H18").Value = 0.930678514466365
$ws.Range("I18").Value = 0.7044726193565225
$ws.Range("M18").Value = 0.4971993356228523

# Row 19
$ws.Range("B19").Value = 1.128042537722479
$ws.Range("C19").Value = 0.3949559508386642
$ws.Range("D19").Value = 0.102169577582444
$ws.Range("E19").Value = 0.1136558544894015
$ws.Range("G19").Value = 1.076989478731434
$ws.Range("H19").Value = 0.9293075964362458
$ws.Range("I19").Value = 0.7034496662119025
$ws.Range("M19").Value = 0.4935639395381912

# Row 20
$ws.Range("B20").Value = 1.168339205485154
$ws.Range("C20").Value = 0.4101140273125452
$ws.Range("D20").Value = 0.1053645689795388
$ws.Range("E20").Value = 0.1160343623439601
$ws.Range("G20").Value = 1.094837019758188
$ws.Range("H20").Value = 0.9355204068186538
$ws.Range("I20").Value = 0.7080904287685001
$ws.Range("M20").Value = 0.5099359105376635

# Row 21
$ws.Range("B21").Value = 1.303944552647067
$ws.Range("C21").Value = 0.4610855759863171
$ws.Range("D21").Value = 0.1161557426894575
$ws.Range("E21").Value = 0.1241310933993915
$ws.Range("G21").Value = 1.155855623711886
$ws.Range("H21").Value = 0.9571608471944444
$ws.Range("I21").Value = 0.7243414211763195
$ws.Range("M21").Value = 0.5651616899913137

# Row 22
$ws.Range("B22").Value = 1.392732141102556
$ws.Range("C22").Value = 0.4944329004465544
$ws.Range("D22").Value = 0.1232513633138126
$ws.Range("E22").Value = 0.1295029940881207
$ws.Range("G22").Value = 1.196532978123599
$ws.Range("H22").Value = 0.9718807359496395
$ws.Range("I22").Value = 0.7354592187850528
$ws.Range("M22").Value = 0.6014215159532341

# Row 23
$ws.Range("B23").Value = 1.345329601055084
$ws.Range("C23").Value = 0.4766314692486731
$ws.Range("D23").Value = 0.1194603076067011
$ws.Range("E23").Value = 0.1266284559023774
$ws.Range("G23").Value = 1.174748829922407
$ws.Range("H23").Value = 0.9639713385763571
$ws.Range("I23").Value = 0.7294795989524516
$ws.Range("M23").Value = 0.5820534646741038

# Row 24
$ws.Range("B24").Value = 1.16612550483427
$ws.Range("C24").Value = 0.4092814653779442
$ws.Range("D24").Value = 0.1051889052604764
$ws.Range("E24").Value = 0.1159033544813113
$ws.Range("G24").Value = 1.09385298878837
$ws.Range("H24").Value = 0.9351763602876133
$ws.Range("I24").Value = 0.707833117377092
$ws.Range("M24").Value = 0.509036029212524

# Row 25
$ws.Range("B25").Value = 0.9735049618369658
$ws.Range("C25").Value = 0.3367633050684731
$ws.Range("D25").Value = 0.08997349261345278
$ws.Range("E25").Value = 0.104667600130945
$ws.Range("G25").Value = 1.009935350906687
$ws.Range("H25").Value = 0.9065561335584107
$ws.Range("I25").Value = 0.6865785576080725
$ws.Range("M25").Value = 0.4309643537454804
